# Actualización automática 2025-05-29 10:10:07
#
# Inserts a new advisor data row ("MACHARE BARCO LISSETTE STEFANIA") right
# above the existing summary/total row on both worksheets, pushing the
# summary row down by one, and updates the "de 1" -> "de 2" counters in the
# summary row of the first sheet (the summary now counts 2 advisors).
# Also widens column B (CLIENTE) from 27 to 33 characters on both sheets.

$wb = $excel.ActiveWorkbook

# ColumnWidth (character units) round-trips through Excel's internal width
# formula with a fixed +0.8333... offset when saved back to OOXML "width".
# Subtract that offset so the persisted <col width="..."/> lands on exactly
# the target value (33).
$colWidthAdjust = 32.166666666666664

# ---------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Columns.Item(2).ColumnWidth = $colWidthAdjust

# Push the current summary row (row 3) down to row 4, inheriting the
# formatting of the row above it (row 2) for the freshly inserted row 3.
$ws1.Rows.Item(3).Insert()

$ws1.Cells.Item(3, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(3, 2).Value = "MACHARE BARCO LISSETTE STEFANIA"
for ($c = 3; $c -le 14; $c++) {
    $ws1.Cells.Item(3, $c).Value = 0
}

# The old summary row (now row 4) counted against 1 advisor ("X de 1");
# now that there are 2 advisors it should read "X de 2".
for ($c = 3; $c -le 14; $c++) {
    $cell = $ws1.Cells.Item(4, $c)
    $cell.Value = $cell.Value2.Replace("de 1", "de 2")
}

# ---------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Columns.Item(2).ColumnWidth = $colWidthAdjust

# Same shift: push the totals row (row 3) down to row 4.
$ws2.Rows.Item(3).Insert()

$ws2.Cells.Item(3, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(3, 2).Value = "MACHARE BARCO LISSETTE STEFANIA"
for ($c = 3; $c -le 6; $c++) {
    $ws2.Cells.Item(3, $c).Value = 0
}
